$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2080.3467
$ws.Range("I15").Value = 2080.3467
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 6241.0401
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -6072.0401

$ws.Range("H103").Value = 853.5714
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 853.5714
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2560.7142
$ws.Range("N103").Value = -3732.7142
$ws.Range("M103").ClearContents()

$ws.Range("H132").Value = 2859318
$ws.Range("I132").Value = 3176246.8
$ws.Range("J132").Value = 6958.857
$ws.Range("K132").Value = 9528740.399999999
$ws.Range("L132").Value = 20876.571
$ws.Range("M132").Value = -9526210.399999999
$ws.Range("N132").Value = -25936.571

$ws.Range("H137").Value = 2861823
$ws.Range("I137").Value = 4550906.5
$ws.Range("J137").Value = 3373.7693
$ws.Range("K137").Value = 13652719.5
$ws.Range("L137").Value = 10121.3079
$ws.Range("M137").Value = -13650169.5
$ws.Range("N137").Value = -15221.3079

$ws.Range("H138").Value = 6017.541
$ws.Range("I138").Value = 2962.76
$ws.Range("J138").Value = 7290.3667
$ws.Range("K138").Value = 8888.280000000001
$ws.Range("L138").Value = 21871.1001
$ws.Range("M138").Value = -3748.280000000001
$ws.Range("N138").Value = -32151.1001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22759.6
$ws.Range("I32").Value = 18942.285
$ws.Range("J32").Value = 31666.666
$ws.Range("K32").Value = 18942.285
$ws.Range("L32").Value = 31666.666
$ws.Range("M32").Value = -18655.285

$ws.Range("H37").Value = 18038
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 18038
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 18038
$ws.Range("N37").Value = -18584

$ws.Range("H110").Value = 3328.1667
$ws.Range("I110").Value = 670.8889
$ws.Range("J110").Value = 11300
$ws.Range("K110").Value = 670.8889
$ws.Range("L110").Value = 11300
$ws.Range("M110").Value = 1374.1111
$ws.Range("N110").Value = -15390

$ws.Range("H122").Value = 4381.6
$ws.Range("I122").Value = 3270.6667
$ws.Range("J122").Value = 4857.7144
$ws.Range("K122").Value = 9812.000100000001
$ws.Range("L122").Value = 14573.1432
$ws.Range("M122").Value = -7362.000100000001
$ws.Range("N122").Value = -19473.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 10007
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10007
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10007
$ws.Range("N15").Value = -10461

$ws.Range("H35").Value = 22905.5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 22905.5
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 22905.5
$ws.Range("N35").Value = -23525.5

$ws.Range("H64").Value = 543
$ws.Range("I64").Value = 469.1
$ws.Range("J64").Value = 599.8461
$ws.Range("K64").Value = 469.1
$ws.Range("L64").Value = 599.8461
$ws.Range("M64").Value = -244.1
$ws.Range("N64").Value = -1049.8461

$ws.Range("H67").Value = 543
$ws.Range("I67").Value = 469.1
$ws.Range("J67").Value = 599.8461
$ws.Range("K67").Value = 469.1
$ws.Range("L67").Value = 599.8461
$ws.Range("M67").Value = 310.9
$ws.Range("N67").Value = -2159.8461

$ws.Range("H82").Value = 22393.455
$ws.Range("I82").Value = 4933
$ws.Range("J82").Value = 28941.125
$ws.Range("K82").Value = 4933
$ws.Range("L82").Value = 28941.125
$ws.Range("M82").Value = -4550
$ws.Range("N82").Value = -29707.125

$ws.Range("H85").Value = 22393.455
$ws.Range("I85").Value = 4933
$ws.Range("J85").Value = 28941.125
$ws.Range("K85").Value = 4933
$ws.Range("L85").Value = 28941.125
$ws.Range("M85").Value = -3607
$ws.Range("N85").Value = -31593.125

$ws.Range("H141").Value = 29700
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 29700
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 29700
$ws.Range("N141").Value = -40060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 4813
$ws.Range("I41").Value = 5100
$ws.Range("J41").Value = 4621.6665
$ws.Range("K41").Value = 5100
$ws.Range("L41").Value = 4621.6665
$ws.Range("M41").Value = -4672
$ws.Range("N41").Value = -5477.6665

$ws.Range("H50").Value = 22092
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 22092
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 22092
$ws.Range("N50").Value = -23342

$ws.Range("H51").Value = 22099
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 22099
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 22099
$ws.Range("N51").Value = -23571

$ws.Range("H59").Value = 22346.834
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 22346.834
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 22346.834
$ws.Range("N59").Value = -24636.834

$ws.Range("H60").Value = 14801.5
$ws.Range("I60").Value = 9500
$ws.Range("J60").Value = 20103
$ws.Range("K60").Value = 9500
$ws.Range("L60").Value = 20103
$ws.Range("M60").Value = -8989
$ws.Range("N60").Value = -21125

$ws.Range("H61").Value = 22099
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 22099
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 22099
$ws.Range("N61").Value = -22795

$ws.Range("H68").Value = 50295
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 50295
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 50295
$ws.Range("N68").Value = -51793

$ws.Range("H71").Value = 50295
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 50295
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 150885
$ws.Range("N71").Value = -158373

$ws.Range("H74").Value = 24608
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 24608
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 24608
$ws.Range("N74").Value = -26356

$ws.Range("H77").Value = 24608
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 24608
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 73824
$ws.Range("N77").Value = -82560

$ws.Range("H122").Value = 2942.5386
$ws.Range("I122").Value = 2622.5
$ws.Range("J122").Value = 3662.625
$ws.Range("K122").Value = 7867.5
$ws.Range("L122").Value = 10987.875
$ws.Range("M122").Value = -5417.5

$ws.Range("H141").Value = 30962.162
$ws.Range("I141").Value = 12500
$ws.Range("J141").Value = 32017.143
$ws.Range("K141").Value = 12500
$ws.Range("L141").Value = 32017.143
$ws.Range("M141").Value = -7320
$ws.Range("N141").Value = -42377.143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 845.62964
$ws.Range("I5").Value = 520.0833
$ws.Range("J5").Value = 3450
$ws.Range("K5").Value = 1560.2499
$ws.Range("L5").Value = 10350
$ws.Range("M5").Value = -1448.2499
$ws.Range("N5").Value = -10574

$ws.Range("H86").Value = 430.83334
$ws.Range("I86").Value = 159.66667
$ws.Range("J86").Value = 702
$ws.Range("K86").Value = 479.00001
$ws.Range("L86").Value = 2106
$ws.Range("M86").Value = 706.99999
$ws.Range("N86").Value = -4478

$ws.Range("H89").Value = 430.83334
$ws.Range("I89").Value = 159.66667
$ws.Range("J89").Value = 702
$ws.Range("K89").Value = 1437.00003
$ws.Range("L89").Value = 6318
$ws.Range("M89").Value = 4490.99997
$ws.Range("N89").Value = -18174

$ws.Range("H122").Value = 8144.4
$ws.Range("I122").Value = 564
$ws.Range("J122").Value = 13198
$ws.Range("K122").Value = 5076
$ws.Range("L122").Value = 118782
$ws.Range("M122").Value = -2626
$ws.Range("N122").Value = -123682

$ws.Range("H135").Value = 845.62964
$ws.Range("I135").Value = 520.0833
$ws.Range("J135").Value = 3450
$ws.Range("K135").Value = 4680.7497
$ws.Range("L135").Value = 31050
$ws.Range("M135").Value = -2145.7497
$ws.Range("N135").Value = -36120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5387.095
$ws.Range("I122").Value = 930
$ws.Range("J122").Value = 7615.643
$ws.Range("K122").Value = 2790
$ws.Range("L122").Value = 22846.929
$ws.Range("M122").Value = -340
$ws.Range("N122").Value = -27746.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 599.5
$ws.Range("I100").Value = 549.4167
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 1098.8334
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -557.8334
$ws.Range("N100").Value = -2882

$ws.Range("H110").Value = 37322
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 37322
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 37322
$ws.Range("N110").Value = -45502

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H119").Value = 25671.143
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 25671.143
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 25671.143
$ws.Range("N119").Value = -35347.143

$ws.Range("H122").Value = 1739.1177
$ws.Range("I122").Value = 1233.2142
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 3699.6426
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -1249.6426
$ws.Range("N122").Value = -17200
